$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark (it currently sits right after
#    "garantias" in the warranty placeholder cell).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Change "       Abono " -> "       Abonos" (drop the trailing space, add
#    the "s") in the table-header cell.
$d.Content.Find.Execute("       Abono ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "       Abonos", 2)

# 3) Word leaves its "_GoBack" bookmark (collapsed) right at the location of
#    the most recent edit - i.e. right after the newly typed "Abonos".
$r = $d.Content
$r.Find.Execute("Abonos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
